# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted above the current row 112
# (Arándano (blue) @ Vega Modelo de Temuco), pushing every subsequent
# record down by one row (112->113, 113->114, ... 122->123, old 123->124).
# The observation that lands on row 123 also has its "Origen" (Región)
# corrected from "Región de La Araucanía" to "Región del Maule".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 112, shifting rows 112:123 down to 113:124.
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly record.
$ws.Range("A112").Value = 10
$ws.Range("B112").Value = "Vega Modelo de Temuco"
$ws.Range("C112").Value = "La Araucanía"
$ws.Range("D112").Value = 44918
$ws.Range("E112").Value = 9
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100101
$ws.Range("H112").Value = "Berries"
$ws.Range("I112").Value = 100101001
$ws.Range("J112").Value = "Arándano (blue)"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 200
$ws.Range("N112").Value = 2000
$ws.Range("O112").Value = 2000
$ws.Range("P112").Value = 2000
$ws.Range("Q112").Value = "$/kilo"
$ws.Range("R112").Value = "Región del Maule"
$ws.Range("S112").Value = 2000
$ws.Range("T112").Value = 1

# The record that shifted from old row 122 down to row 123 also gets its
# Origen corrected.
$ws.Range("R123").Value = "Región del Maule"
